$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.08995199074074074
$ws.Range("B2").Value = 35.05374861111112
$ws.Range("B3").Value = 1804.831415305278
$ws.Range("B5").Value = 3.902
$ws.Range("B7").Value = 9
$ws.Range("B8").Value = 67.44375803705395
$ws.Range("B9").Value = 26.76054045377623
$ws.Range("B10").Value = 88
$ws.Range("B11").Value = "Eco mode`n100.00%"
$ws.Range("B13").Value = -839.3480293371658
$ws.Range("B14").Value = 19.54605749888889
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 1.071382309322508
$ws.Range("B17").Value = 3.024
$ws.Range("B18").Value = 0.3149999999999999
$ws.Range("B32").Value = 1.927956173611111
$ws.Range("B33").Value = 0.00000006891563268030396
$ws.Range("B34").Value = 18.10137851034557
$ws.Range("B35").Value = 3.625671045887399
$ws.Range("B36").Value = 1.754835576897138
$ws.Range("B37").Value = 2.947206560738083
$ws.Range("B38").Value = 41.9757749062559
$ws.Range("B39").Value = 29.21983328387602
